# Update TPM-derived ligand/receptor specificity and edge-weight values on Sheet1.
# Root inputs that changed are the ligand average/total expression for the
# "ECs" sending cluster (G2:H4) and the receptor average/total expression for
# the "ECs" target cluster (M2/N2, M5/N5, M8/N8). Every other touched column
# (I, J, O, P = derived specificities; Q, R, S, T = edge weights/specificities)
# is a downstream recomputation of those roots. Values below are written
# directly with the final recomputed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Sending=ECs, Target=ECs)
$ws.Range("G2").Value = 0.5979736666666666
$ws.Range("H2").Value = 1.793921
$ws.Range("I2").Value = 0.03342655292740804
$ws.Range("J2").Value = 0.03342655292740804
$ws.Range("M2").Value = 0.8213140000000001
$ws.Range("N2").Value = 2.463942
$ws.Range("O2").Value = 0.06824749762056036
$ws.Range("P2").Value = 0.06824749762056037
$ws.Range("Q2").Value = 0.4911241440646667
$ws.Range("R2").Value = 4.420117296582
$ws.Range("S2").Value = 0.002281278591376815
$ws.Range("T2").Value = 0.002281278591376816

# Row 3 (Sending=ECs, Target=FAPs)
$ws.Range("G3").Value = 0.5979736666666666
$ws.Range("H3").Value = 1.793921
$ws.Range("I3").Value = 0.03342655292740804
$ws.Range("J3").Value = 0.03342655292740804
$ws.Range("O3").Value = 0.2017018900182306
$ws.Range("P3").Value = 0.2017018900182306
$ws.Range("Q3").Value = 1.451491579107889
$ws.Range("R3").Value = 13.063424211971
$ws.Range("S3").Value = 0.006742198902252622
$ws.Range("T3").Value = 0.006742198902252622

# Row 4 (Sending=ECs, Target=MuSCs)
$ws.Range("G4").Value = 0.5979736666666666
$ws.Range("H4").Value = 1.793921
$ws.Range("I4").Value = 0.03342655292740804
$ws.Range("J4").Value = 0.03342655292740804
$ws.Range("O4").Value = 0.7300506123612091
$ws.Range("P4").Value = 0.7300506123612091
$ws.Range("Q4").Value = 5.253606280382778
$ws.Range("R4").Value = 47.282456523445
$ws.Range("S4").Value = 0.02440307543377861
$ws.Range("T4").Value = 0.02440307543377861

# Row 5 (Sending=FAPs, Target=ECs)
$ws.Range("I5").Value = 0.8874158839838097
$ws.Range("J5").Value = 0.8874158839838097
$ws.Range("M5").Value = 0.8213140000000001
$ws.Range("N5").Value = 2.463942
$ws.Range("O5").Value = 0.06824749762056036
$ws.Range("P5").Value = 0.06824749762056037
$ws.Range("Q5").Value = 13.03847774544467
$ws.Range("R5").Value = 117.346299709002
$ws.Range("S5").Value = 0.06056391343063251
$ws.Range("T5").Value = 0.06056391343063253

# Row 6 (Sending=FAPs, Target=FAPs)
$ws.Range("I6").Value = 0.8874158839838097
$ws.Range("J6").Value = 0.8874158839838097
$ws.Range("O6").Value = 0.2017018900182306
$ws.Range("P6").Value = 0.2017018900182306
$ws.Range("S6").Value = 0.1789934610317333
$ws.Range("T6").Value = 0.1789934610317333

# Row 7 (Sending=FAPs, Target=MuSCs)
$ws.Range("I7").Value = 0.8874158839838097
$ws.Range("J7").Value = 0.8874158839838097
$ws.Range("O7").Value = 0.7300506123612091
$ws.Range("P7").Value = 0.7300506123612091
$ws.Range("S7").Value = 0.6478585095214439
$ws.Range("T7").Value = 0.6478585095214439

# Row 8 (Sending=MuSCs, Target=ECs)
$ws.Range("I8").Value = 0.07915756308878232
$ws.Range("J8").Value = 0.07915756308878232
$ws.Range("M8").Value = 0.8213140000000001
$ws.Range("N8").Value = 2.463942
$ws.Range("O8").Value = 0.06824749762056036
$ws.Range("P8").Value = 0.06824749762056037
$ws.Range("Q8").Value = 1.163033188096
$ws.Range("R8").Value = 10.467298692864
$ws.Range("S8").Value = 0.005402305598551028
$ws.Range("T8").Value = 0.005402305598551029

# Row 9 (Sending=MuSCs, Target=FAPs)
$ws.Range("I9").Value = 0.07915756308878232
$ws.Range("J9").Value = 0.07915756308878232
$ws.Range("O9").Value = 0.2017018900182306
$ws.Range("P9").Value = 0.2017018900182306
$ws.Range("Q9").Value = 3.437283422421334
$ws.Range("S9").Value = 0.01596623008424472
$ws.Range("T9").Value = 0.01596623008424472

# Row 10 (Sending=MuSCs, Target=MuSCs)
$ws.Range("I10").Value = 0.07915756308878232
$ws.Range("J10").Value = 0.07915756308878232
$ws.Range("O10").Value = 0.7300506123612091
$ws.Range("P10").Value = 0.7300506123612091
$ws.Range("S10").Value = 0.05778902740598657
$ws.Range("T10").Value = 0.05778902740598657

Write-Host "Updated TPM-derived values on $($ws.Name)"
